$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(10, 3).Value = "'314"
$ws.Cells.Item(10, 4).Value = "'1080587.67"
$ws.Cells.Item(20, 3).Value = "'151"
$ws.Cells.Item(20, 4).Value = "'386569.00"
$ws.Cells.Item(21, 3).Value = "'308"
$ws.Cells.Item(21, 4).Value = "'1014217.92"
$ws.Cells.Item(22, 3).Value = "'149"
$ws.Cells.Item(22, 4).Value = "'403137.26"
$ws.Cells.Item(24, 3).Value = "'7"
$ws.Cells.Item(24, 4).Value = "'22700.00"
$ws.Cells.Item(30, 3).Value = "'427"
$ws.Cells.Item(30, 4).Value = "'1478959.16"
$ws.Cells.Item(32, 3).Value = "'329"
$ws.Cells.Item(32, 4).Value = "'985882.96"
$ws.Cells.Item(40, 3).Value = "'74"
$ws.Cells.Item(40, 4).Value = "'193746.00"
$ws.Cells.Item(41, 3).Value = "'51"
$ws.Cells.Item(41, 4).Value = "'236797.92"
$ws.Cells.Item(42, 3).Value = "'86"
$ws.Cells.Item(42, 4).Value = "'381157.99"
$ws.Cells.Item(43, 3).Value = "'5"
$ws.Cells.Item(43, 4).Value = "'23132.00"
$ws.Cells.Item(57, 3).Value = "'129"
$ws.Cells.Item(57, 4).Value = "'277500.00"
$ws.Cells.Item(58, 3).Value = "'221"
$ws.Cells.Item(58, 4).Value = "'598500.00"
$ws.Cells.Item(59, 3).Value = "'157"
$ws.Cells.Item(59, 4).Value = "'382439.00"
$ws.Cells.Item(60, 3).Value = "'6"
$ws.Cells.Item(60, 4).Value = "'13500.00"
$ws.Cells.Item(68, 3).Value = "'339"
$ws.Cells.Item(68, 4).Value = "'831135.70"
$ws.Cells.Item(70, 3).Value = "'830"
$ws.Cells.Item(70, 4).Value = "'2670705.34"
$ws.Cells.Item(71, 3).Value = "'474"
$ws.Cells.Item(71, 4).Value = "'1434386.03"
$ws.Cells.Item(73, 3).Value = "'30"
$ws.Cells.Item(73, 4).Value = "'103736.09"
$ws.Cells.Item(74, 3).Value = "'392"
$ws.Cells.Item(74, 4).Value = "'894342.96"
$ws.Cells.Item(77, 3).Value = "'554"
$ws.Cells.Item(77, 4).Value = "'1431857.16"
